$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Add the hidden "DropdownOptions" sheet right after Sheet1 ---------
$dropdown = $wb.Worksheets.Add($null, $ws)
$dropdown.Name = "DropdownOptions"

$options = @("0% - 10%", "11% - 25%", "26% - 50%", "51% - 75%", "76% - 90%", "91% - 99%")
for ($i = 0; $i -lt $options.Length; $i++) {
  $dropdown.Cells.Item($i + 1, 1).Value = $options[$i]
}
# "100%" would auto-coerce to the number 1 with a percent format unless the
# cell is pre-formatted as text; set it as text, write it, then drop the
# number-format override so the cell ends up as a plain string again.
$dropdown.Cells.Item(7, 1).NumberFormat = "@"
$dropdown.Cells.Item(7, 1).Value = "100%"
$dropdown.Cells.Item(7, 1).ClearFormats()

$dropdown.Visible = $false

# --- 2. New column AU: "Status as of July 4, 2025" ------------------------
$ws.Cells.Item(1, 47).Value = "Status as of July 4, 2025"

# --- 3. Drop the stray "blank" inline-string cells scattered in rows 2-7 --
$ws.Range("P2:S2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AM2").ClearContents()
$ws.Range("AP2:AR2").ClearContents()

$ws.Range("P3:Y3").ClearContents()
$ws.Range("AM3").ClearContents()
$ws.Range("AO3:AR3").ClearContents()

$ws.Range("P4:Q4").ClearContents()
$ws.Range("S4").ClearContents()
$ws.Range("Z4").ClearContents()
$ws.Range("AM4").ClearContents()
$ws.Range("AO4:AR4").ClearContents()

$ws.Range("P5:Q5").ClearContents()
$ws.Range("S5").ClearContents()
$ws.Range("Z5").ClearContents()
$ws.Range("AM5").ClearContents()
$ws.Range("AO5:AR5").ClearContents()

$ws.Range("P6:Q6").ClearContents()
$ws.Range("S6").ClearContents()
$ws.Range("Z6").ClearContents()
$ws.Range("AM6").ClearContents()
$ws.Range("AO6:AR6").ClearContents()

$ws.Range("P7:X7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AM7").ClearContents()
$ws.Range("AO7:AR7").ClearContents()

# --- 4. Re-apply the date/time number format on the bid-schedule columns --
$ws.Range("T2:W2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("T4:X4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("T5:X5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("T6:X6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- 5. Dropdown validation on the new column, driven by the hidden sheet -
$validationRange = $ws.Range("AU2:AU7")
$validationRange.Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$validationRange.Validation.IgnoreBlank = $true
$validationRange.Validation.InCellDropdown = $true
$validationRange.Validation.ShowInput = $false
$validationRange.Validation.ShowError = $false
